# fix(data): add real ids
# Replace the sequential numeric IDs in column A with the real room id
# (a UUID) that every row shares.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newId = "0195e7a9-181c-728b-bc78-0a90c0cd8d0f"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 6 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newId
}
